$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("F2").Value = 1.87
$ws.Range("G2").Value = 1.88
$ws.Range("N2").Value = 3.45
$ws.Range("R2").Value = 1.3
$ws.Range("W2").Value = 2.12
$ws.Range("Y2").Value = 15.5
$ws.Range("AB2").Value = 8
$ws.Range("AD2").Value = 19.5
$ws.Range("AG2").Value = 10
$ws.Range("AI2").Value = 85
$ws.Range("AO2").Value = 95

# Row 3 updates
$ws.Range("G3").Value = 5.7
$ws.Range("H3").Value = 1.81
$ws.Range("I3").Value = 1.82
$ws.Range("J3").Value = 3.7
$ws.Range("K3").Value = 3.75
$ws.Range("P3").Value = 1.78
$ws.Range("Q3").Value = 2.22
$ws.Range("R3").Value = 1.29
$ws.Range("W3").Value = 1.21
$ws.Range("AB3").Value = 16
$ws.Range("AO3").Value = 14.5

# Row 4 updates
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 3
$ws.Range("AF4").Value = 46
$ws.Range("AG4").Value = 25
$ws.Range("AJ4").Value = 190
$ws.Range("AM4").Value = 210
$ws.Range("AN4").Value = 180
$ws.Range("AO4").Value = 17

# Row 5 updates
$ws.Range("H5").Value = 1.77
$ws.Range("I5").Value = 1.78
$ws.Range("L5").Value = 1.43
$ws.Range("M5").Value = 1.08
$ws.Range("T5").Value = 1.96
$ws.Range("V5").Value = 2.28
$ws.Range("AH5").Value = 21

# Row 6 updates
$ws.Range("F6").Value = 2.8
$ws.Range("H6").Value = 2.68
$ws.Range("L6").Value = 1.42
$ws.Range("M6").Value = 1.08
$ws.Range("S6").Value = 3.85
$ws.Range("U6").Value = 2.14
$ws.Range("AA6").Value = 55
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 32
$ws.Range("AI6").Value = 48
$ws.Range("AK6").Value = 36
$ws.Range("AO6").Value = 30

# Row 7 updates
$ws.Range("L7").Value = 1.4
$ws.Range("P7").Value = 1.97
$ws.Range("U7").Value = 2.24
$ws.Range("AG7").Value = 14.5
$ws.Range("AO7").Value = 19.5

# Row 8 updates
$ws.Range("G8").Value = 1.65
$ws.Range("H8").Value = 6.2
$ws.Range("I8").Value = 7.2
$ws.Range("K8").Value = 4.4
$ws.Range("L8").Value = 1.32
$ws.Range("O8").Value = 1.23
$ws.Range("Q8").Value = 1.75
$ws.Range("S8").Value = 2.66
$ws.Range("T8").Value = 1.78
$ws.Range("U8").Value = 2.14
$ws.Range("V8").Value = 1.17
$ws.Range("W8").Value = 2.52
$ws.Range("X8").Value = 24
$ws.Range("Z8").Value = 70
$ws.Range("AA8").Value = 1000
$ws.Range("AC8").Value = 9.6
$ws.Range("AD8").Value = 24
$ws.Range("AE8").Value = 80
$ws.Range("AF8").Value = 10
$ws.Range("AH8").Value = 18.5
$ws.Range("AI8").Value = 80
